# Apply the vocabulary.xlsx update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the value to be stored as literal text (columns like Y hold
    # dates written as plain strings, e.g. "2024-01-08" -- without this,
    # Excel's autoconvert would turn it into a date serial number).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Row 50: extend the G50 list, update Y50 date ---
$ws.Cells.Item(50, 7).Value = "vocab:1000,vocab:1038,vocab:1034,vocab:1028,vocab:1064,vocab:1057,vocab:1252,vocab:1244,vocab:1263,vocab:1259"
Set-TextValue $ws.Cells.Item(50, 25) "2024-01-08"

# --- Rows 279-282: rename terms, update dates ---
$ws.Cells.Item(279, 2).Value = "urinary samples information"
Set-TextValue $ws.Cells.Item(279, 25) "2024-01-08"

$ws.Cells.Item(280, 2).Value = "creatinine"
Set-TextValue $ws.Cells.Item(280, 25) "2024-01-08"

$ws.Cells.Item(281, 2).Value = "specific gravity"
Set-TextValue $ws.Cells.Item(281, 25) "2024-01-08"

$ws.Cells.Item(282, 2).Value = "osmolarity"
Set-TextValue $ws.Cells.Item(282, 25) "2024-01-08"

# --- New rows 283-307 ---
# columns: A=id, B=term, G=parent(s), Y=date; all other columns blank
$newRows = @(
    @{ Row = 283; A = "vocab:1263"; B = "blood samples information"; G = "" },
    @{ Row = 284; A = "vocab:1264"; B = "gravimetric"; G = "vocab:1263" },
    @{ Row = 285; A = "vocab:1265"; B = "enzymatic"; G = "vocab:1263" },
    @{ Row = 286; A = "vocab:1266"; B = "none"; G = "vocab:1263,vocab:1259" },
    @{ Row = 287; A = "vocab:1267"; B = "unit of measure"; G = "" },
    @{ Row = 288; A = "vocab:1268"; B = "mL"; G = "vocab:1267" },
    @{ Row = 289; A = "vocab:1269"; B = "µL"; G = "vocab:1267" },
    @{ Row = 290; A = "vocab:1270"; B = "g"; G = "vocab:1267" },
    @{ Row = 291; A = "vocab:1271"; B = "mg"; G = "vocab:1267" },
    @{ Row = 292; A = "vocab:1272"; B = "µg"; G = "vocab:1267" },
    @{ Row = 293; A = "vocab:1273"; B = "cm"; G = "vocab:1267" },
    @{ Row = 294; A = "vocab:1274"; B = "sampling container material"; G = "" },
    @{ Row = 295; A = "vocab:1275"; B = "PP"; G = "vocab:1274" },
    @{ Row = 296; A = "vocab:1276"; B = "glass"; G = "vocab:1274" },
    @{ Row = 297; A = "vocab:1277"; B = "PET"; G = "vocab:1274" },
    @{ Row = 298; A = "vocab:1278"; B = "quality assurance/qualtiy control method"; G = "" },
    @{ Row = 299; A = "vocab:1279"; B = "standard operating procedure"; G = "vocab:1278" },
    @{ Row = 300; A = "vocab:1280"; B = "trained fieldworkers"; G = "vocab:1278" },
    @{ Row = 301; A = "vocab:1281"; B = "control of background contamination in the sampling material"; G = "vocab:1278" },
    @{ Row = 302; A = "vocab:1282"; B = "controlf of the transprot conditions"; G = "vocab:1278" },
    @{ Row = 303; A = "vocab:1283"; B = "control of background contamination in the conservation material"; G = "vocab:1278" },
    @{ Row = 304; A = "vocab:1284"; B = "identifciation/traceability of the samples"; G = "vocab:1278" },
    @{ Row = 305; A = "vocab:1285"; B = "criteria for acceptation/acceptation of the samples"; G = "vocab:1278" },
    @{ Row = 306; A = "vocab:1286"; B = "collection of field blanks"; G = "vocab:1278" },
    @{ Row = 307; A = "vocab:1287"; B = "settings the conditions for sample storage"; G = "vocab:1278" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    if ($r.G -ne "") {
        $ws.Cells.Item($r.Row, 7).Value = $r.G
    }
    Set-TextValue $ws.Cells.Item($r.Row, 25) "2024-01-08"
}
